$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = "System"
$ws.Cells.Item(7, 6).Value = "2025-03-27 17:47:33"
$ws.Cells.Item(7, 7).Value = 0

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 18
$ws.Cells.Item(8, 3).Value = 7
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = "System"
$ws.Cells.Item(8, 6).Value = "2025-03-27 17:50:33"
$ws.Cells.Item(8, 7).Value = 0
